# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# Update OFF (offense) target depth data - Road ("R") row
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 209
$wsOff.Range("C3").Value = 151
$wsOff.Range("D3").Value = 50
$wsOff.Range("E3").Value = 24
$wsOff.Range("G3").Value = 4

# Update DEF (defense) target depth data - Road ("R") row
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 138
$wsDef.Range("C3").Value = 92
$wsDef.Range("D3").Value = 31
$wsDef.Range("E3").Value = 15
$wsDef.Range("F3").Value = 3
